# "Riesgos de un proyecto de tics pdf"
#
# The "Medidas a tomar" (measures to take) column (D) had its four
# descriptive texts rewritten with embedded line breaks (and one typo
# fix: "perioso" -> "periodo"), the affected cells were switched to a
# wrapped-text style with taller rows, column D was narrowed, and the
# view was re-zoomed with a new selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated risk-mitigation texts (column D) ---------------------------
# NB: written in this particular order so the workbook's shared-string
# table is rebuilt in the same sequence as the target file.
$ws.Range("D6").Value = "Dependiendo el tipo de proyecto, podemos pedir mas financiamiento, `no bien pedir apoyo por otros medios como donaciones o incluso dejar`n que la empresa sea absorvida por una empresa mas grande"
$ws.Range("D3").Value = "Dedicar equipos que puedan brindar apoyo y contenido para aquellos `nque tienen alguna duda, con la aplicacion, motor, requerimientos"
$ws.Range("D8").Value = "Hacer que todo el equipo este conciente de los cambios, evaluar`ncuales son los mas importantes y empezar a trabajar en ellos antes `nde continuar con otros requerimientos  "
$ws.Range("D4").Value = "Realizar testings al finalizar un requerimiento que se puede considerar`n grande, y no al final del sprint en el periodo de lanzamiento a produccion  "

# --- Wrap text + taller rows for the four updated cells ------------------
$ws.Range("D3").WrapText = $true
$ws.Range("D4").WrapText = $true
$ws.Range("D6").WrapText = $true
$ws.Range("D8").WrapText = $true

$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 45

# --- Narrower column D (was 177 characters wide) --------------------------
$ws.Range("D1").EntireColumn.ColumnWidth = 64.3

# --- View: zoom in and move the selection -------------------------------
$excel.ActiveWindow.Zoom = 160
[void]$ws.Range("G4").Select()
